$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the three paragraphs following "BADAN PUSAT STATISTIK" that hold
#    the (now unused) small address / fax block:
#      - empty sz=10 spacer paragraph
#      - "Jalan dr. Sutomo No. 6-8 Jakarta 10710, ..."
#      - "Fax. (021) 3857046, Homepage: ..."
# ---------------------------------------------------------------------------
$badanIdx = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text -like '*BADAN PUSAT STATISTIK*') {
        $badanIdx = $i
        break
    }
}

if ($badanIdx -gt 0) {
    $firstToRemove = $d.Paragraphs.Item($badanIdx + 1)
    $lastToRemove  = $d.Paragraphs.Item($badanIdx + 3)
    $killRange = $d.Range($firstToRemove.Range.Start, $lastToRemove.Range.End)
    $killRange.Delete()
}

# ---------------------------------------------------------------------------
# 2) Drop the centered alignment on the (now) blank paragraph right after the
#    bookmark end (it no longer should be forced to center).
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like '*BADAN PUSAT STATISTIK*') {
        $p.Next().Alignment = 0
        break
    }
}

# ---------------------------------------------------------------------------
# 3) Remove the justified ("both") alignment on the table-cell paragraph that
#    holds the "${nama}" placeholder.
# ---------------------------------------------------------------------------
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like '*${nama}*' -and $p.Range.Text -notlike '*nama_pengirim*') {
        $p.Alignment = 0
        break
    }
}
